# Activity_Table.xlsx edit
# "them hoat dong truoc cho cac code module"
#
# The "Thiet ke giao dien" (F) and "Xay dung giao dien" (G) activities are
# merged into a single activity "Thiet ke va xay dung giao dien" (F) with a
# duration of 3 days. The former "Hoan chinh giao dien" activity moves up to
# take slot G. All the "Code module ..." activities (previously H..N with
# predecessor H) move up one slot (now H..N) and now depend on "D, G" instead
# of just "H". "Kiem thu cac module" moves up to slot O and its predecessor
# list is updated to "H, I, J, K, L, M, N". "Hoan thien san pham" moves up to
# slot P with predecessor "O", "Bao cao" moves up to slot Q with predecessor
# "P", and the final row (previously R, "Bao cao") is cleared out / removed
# from the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 (ID "F"): merge "Thiet ke giao dien" + "Xay dung giao dien" ---
$ws.Cells.Item(7, 2).Value = "Thiết kế và xây dựng giao diện"
$ws.Cells.Item(7, 3).Value = 3
$ws.Cells.Item(7, 4).Value = "-"

# --- Row 8 (ID "G"): becomes "Hoan chinh giao dien" ---
$ws.Cells.Item(8, 2).Value = "Hoàn chỉnh giao diện"
$ws.Cells.Item(8, 3).Value = 1
$ws.Cells.Item(8, 4).Value = "F"

# --- Rows 9-15 (ID "H".."N"): former code-module rows shift up one slot,
#     predecessor becomes "D, G" ---
$ws.Cells.Item(9, 2).Value = "Code module quản lý khách hàng"
$ws.Cells.Item(9, 3).Value = 2
$ws.Cells.Item(9, 4).Value = "D, G"

$ws.Cells.Item(10, 2).Value = "Code module quản lý nhân viên"
$ws.Cells.Item(10, 3).Value = 2
$ws.Cells.Item(10, 4).Value = "D, G"

$ws.Cells.Item(11, 2).Value = "Code module quản lý phòng"
$ws.Cells.Item(11, 3).Value = 2
$ws.Cells.Item(11, 4).Value = "D, G"

$ws.Cells.Item(12, 2).Value = "Code module quản lý dịch vụ"
$ws.Cells.Item(12, 3).Value = 2
$ws.Cells.Item(12, 4).Value = "D, G"

$ws.Cells.Item(13, 2).Value = "Code module quản lý book"
$ws.Cells.Item(13, 3).Value = 2
$ws.Cells.Item(13, 4).Value = "D, G"

$ws.Cells.Item(14, 2).Value = "Code module quản lý hoá đơn"
$ws.Cells.Item(14, 3).Value = 2
$ws.Cells.Item(14, 4).Value = "D, G"

$ws.Cells.Item(15, 2).Value = "Code module quản lý tài khoản"
$ws.Cells.Item(15, 3).Value = 2
$ws.Cells.Item(15, 4).Value = "D, G"

# --- Row 16 (ID "O"): becomes "Kiem thu cac module" ---
$ws.Cells.Item(16, 2).Value = "Kiểm thử các module"
$ws.Cells.Item(16, 3).Value = 1
$ws.Cells.Item(16, 4).Value = "H, I, J, K, L, M, N"

# --- Row 17 (ID "P"): becomes "Hoan thien san pham" ---
$ws.Cells.Item(17, 2).Value = "Hoàn thiện sản phẩm"
$ws.Cells.Item(17, 3).Value = 1
$ws.Cells.Item(17, 4).Value = "O"

# --- Row 18 (ID "Q"): becomes "Bao cao" ---
$ws.Cells.Item(18, 2).Value = "Báo cáo"
$ws.Cells.Item(18, 3).Value = 1
$ws.Cells.Item(18, 4).Value = "P"

# --- Row 19 (ID "R"): activity removed, clear the row content ---
$ws.Cells.Item(19, 1).ClearContents()
$ws.Cells.Item(19, 2).ClearContents()
$ws.Cells.Item(19, 3).ClearContents()
$ws.Cells.Item(19, 4).ClearContents()

# A19 keeps the "Good" look (green font, centered) but without the green
# highlight fill, matching the new cellXf added for the now-empty ID cell.
$cellA19 = $ws.Cells.Item(19, 1)
$ws.Cells.Item(3, 1).Copy()
$cellA19.PasteSpecial(-4122)
$cellA19.Interior.Pattern = -4142
$excel.CutCopyMode = $false

# --- Shrink the table / autofilter range to exclude the now-empty row 19 ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A2:D18"))

# --- Update the active selection to reflect where the author left off ---
$ws.Range("F19").Select()

$wb.Save()
